# docs_virtocommerce_org-sitemap.xlsx edits
# Commit: "Added Scalability, Feature flags, Import products to catalog, Catalog publishing article"
#
# The sheet tracks documentation pages being ported to a new site. A handful
# of rows that referenced the old "lesson1/lesson2" / ad-hoc "copy content"
# placeholder text are being turned into concrete "OnReview" rows pointing at
# the final article URLs (and the redirect owner/notes columns get filled in
# to match the other OnReview rows). One brand-new row (101) gets populated
# the same way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 60: getting-started/lessons/lesson1 -> OnReview
$ws.Range("B60").Value = "platform/developer-guide/Getting-Started/Installation-Guide/windows/"
$ws.Range("C60").Value = "OnReview"

# Row 61: getting-started/lessons/lesson2 -> Import products to catalog article, OnReview
$ws.Range("B61").Value = "platform/developer-guide/Tutorials-and-How-tos/Tutorials/import-products-to-catalog/"
$ws.Range("C61").Value = "OnReview"

# Row 83: catalog-personalization settings article, OnReview
$ws.Range("B83").Value = "platform/user-guide/catalog-personalization/settings/"
$ws.Range("C83").Value = "OnReview"

# Row 84: Catalog publishing overview article, OnReview
$ws.Range("B84").Value = "platform/user-guide/catalog-publishing/overview/"
$ws.Range("C84").Value = "OnReview"

# Row 90: override-rounding-policy how-to article, OnReview
$ws.Range("B90").Value = "platform/developer-guide/Tutorials-and-How-tos/How-tos/overriding-rounding-policy/"
$ws.Range("C90").Value = "OnReview"

# Row 92: contacts filtering-options article, OnReview
$ws.Range("B92").Value = "platform/user-guide/contacts/filtering-options/"
$ws.Range("C92").Value = "OnReview"

# Row 101: new dynamic-associations overview article, OnReview, owner + note
$ws.Range("B101").Value = "platform/user-guide/marketing/dynamic-associations-overview/"
$ws.Range("C101").Value = "OnReview"
$ws.Range("D101").Value = "Maria"
$ws.Range("E101").Value = "Add redirect when  done"

# Restore the author's last selection/view position on the sheet.
$ws.Activate()
$ws.Range("C61").Select()
